# Add a new LeetCode entry row ("Find the Maximum Length of Valid Subsequence 1")
# to the bottom of the tracking sheet, and move the selection to reflect where
# the author was working (D34) as described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new row (row 47) with Id / Title / Type / Steps.
$ws.Range("A47").Value = 3201
$ws.Range("B47").Value = "Find the Maximum Length of Valid Subsequent 1"
$ws.Range("C47").Value = "Even/Odd, Dynamic Programming"
$ws.Range("D47").Value = "(E+E)/2 = E; (O+O)/2 = E | (E+O)/2=O; (O+E)/2=O. Find the longest subsequence of evens, odds, e-o or o-e. Have a dp to keep track and compare foreach element"

# Update the active selection to match where the author ended up (D34).
[void]$ws.Range("D34").Select()
